# Update Leve profit-tracking sheets with refreshed Market Board price data.
# Each sheet's table (Table_<SheetName>) stores current average prices (H:J)
# and derived Leve price/profit columns (K:N); this applies the latest refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 14162.682
$ws.Range("I18").Value = 9059.9
$ws.Range("J18").Value = 18415
$ws.Range("K18").Value = 9059.9
$ws.Range("L18").Value = 18415
$ws.Range("M18").Value = -8775.9
$ws.Range("N18").Value = -18983
$ws.Range("H19").Value = 1479.8846
$ws.Range("I19").Value = 1278.75
$ws.Range("J19").Value = 1652.2858
$ws.Range("K19").Value = 1278.75
$ws.Range("L19").Value = 1652.2858
$ws.Range("M19").Value = -1103.75
$ws.Range("N19").Value = -2002.2858
$ws.Range("H51").Value = 4455.3335
$ws.Range("J51").Value = 4637.375
$ws.Range("L51").Value = 4637.375
$ws.Range("N51").Value = -5605.375
$ws.Range("H64").Value = 2950
$ws.Range("I64").Value = 3000
$ws.Range("J64").Value = 2900
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 2900
$ws.Range("M64").Value = -2752
$ws.Range("N64").Value = -3396
$ws.Range("H67").Value = 2950
$ws.Range("I67").Value = 3000
$ws.Range("J67").Value = 2900
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 2900
$ws.Range("M67").Value = -2142
$ws.Range("N67").Value = -4616
$ws.Range("H70").Value = 18380.908
$ws.Range("I70").Value = 13551
$ws.Range("J70").Value = 19454.223
$ws.Range("K70").Value = 40653
$ws.Range("L70").Value = 58362.66900000001
$ws.Range("M70").Value = -40383
$ws.Range("N70").Value = -58902.66900000001
$ws.Range("H73").Value = 18380.908
$ws.Range("I73").Value = 13551
$ws.Range("J73").Value = 19454.223
$ws.Range("K73").Value = 40653
$ws.Range("L73").Value = 58362.66900000001
$ws.Range("M73").Value = -39717
$ws.Range("N73").Value = -60234.66900000001
$ws.Range("H103").Value = 822.2222
$ws.Range("I103").Value = 866.6667
$ws.Range("J103").Value = 800
$ws.Range("K103").Value = 2600.0001
$ws.Range("L103").Value = 2400
$ws.Range("M103").Value = -2014.0001
$ws.Range("N103").Value = -3572
$ws.Range("H135").Value = 733
$ws.Range("I135").Value = 708
$ws.Range("K135").Value = 6372
$ws.Range("M135").Value = -3837
$ws.Range("H137").Value = 36402.863
$ws.Range("J137").Value = 40526.04
$ws.Range("L137").Value = 121578.12
$ws.Range("N137").Value = -126678.12

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3447.3618
$ws.Range("I32").Value = 2662.4119
$ws.Range("J32").Value = 5500.3076
$ws.Range("K32").Value = 2662.4119
$ws.Range("L32").Value = 5500.3076
$ws.Range("M32").Value = -2375.4119
$ws.Range("N32").Value = -6074.3076
$ws.Range("H61").Value = 3422.9546
$ws.Range("I61").Value = 2463.6428
$ws.Range("K61").Value = 2463.6428
$ws.Range("M61").Value = -2251.6428
$ws.Range("H74").Value = 1992.6
$ws.Range("J74").Value = 3928.5
$ws.Range("L74").Value = 3928.5
$ws.Range("N74").Value = -5676.5
$ws.Range("H77").Value = 1992.6
$ws.Range("J77").Value = 3928.5
$ws.Range("L77").Value = 19642.5
$ws.Range("N77").Value = -28378.5
$ws.Range("H136").Value = 3422.9546
$ws.Range("I136").Value = 2463.6428
$ws.Range("K136").Value = 7390.928400000001
$ws.Range("M136").Value = -4840.928400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1958.3889
$ws.Range("I20").Value = 1926.9231
$ws.Range("J20").Value = 2040.2
$ws.Range("K20").Value = 1926.9231
$ws.Range("L20").Value = 2040.2
$ws.Range("M20").Value = -1679.9231
$ws.Range("N20").Value = -2534.2
$ws.Range("H105").Value = 2583.3333
$ws.Range("I105").Value = 2643.75
$ws.Range("K105").Value = 2643.75
$ws.Range("M105").Value = -896.75
$ws.Range("H134").Value = 10075.923
$ws.Range("I134").Value = 11189.728
$ws.Range("K134").Value = 33569.18399999999
$ws.Range("M134").Value = -31034.18399999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1366.9714
$ws.Range("I31").Value = 869.1667
$ws.Range("J31").Value = 1894.0588
$ws.Range("K31").Value = 869.1667
$ws.Range("L31").Value = 1894.0588
$ws.Range("M31").Value = -574.1667
$ws.Range("N31").Value = -2484.0588
$ws.Range("H34").Value = 1366.9714
$ws.Range("I34").Value = 869.1667
$ws.Range("J34").Value = 1894.0588
$ws.Range("K34").Value = 869.1667
$ws.Range("L34").Value = 1894.0588
$ws.Range("M34").Value = -667.1667
$ws.Range("N34").Value = -2298.0588
$ws.Range("H105").Value = 1566.6666
$ws.Range("I105").Value = 1900
$ws.Range("K105").Value = 1900
$ws.Range("M105").Value = -153
$ws.Range("H134").Value = 2607
$ws.Range("I134").Value = 1581.4
$ws.Range("K134").Value = 4744.200000000001
$ws.Range("M134").Value = -2209.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 126887.125
$ws.Range("I113").Value = 505001.5
$ws.Range("J113").Value = 849
$ws.Range("K113").Value = 1515004.5
$ws.Range("L113").Value = 2547
$ws.Range("M113").Value = -1512834.5
$ws.Range("N113").Value = -6887
$ws.Range("H131").Value = 6183801
$ws.Range("J131").Value = 11507.922
$ws.Range("L131").Value = 34523.766
$ws.Range("N131").Value = -44603.766
$ws.Range("H134").Value = 5158.8887
$ws.Range("I134").Value = 5476.6665
$ws.Range("K134").Value = 16429.9995
$ws.Range("M134").Value = -11359.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 19665.666
$ws.Range("J34").Value = 19665.666
$ws.Range("L34").Value = 19665.666
$ws.Range("N34").Value = -20201.666
$ws.Range("H70").Value = 5399.8
$ws.Range("J70").Value = 4999.75
$ws.Range("L70").Value = 4999.75
$ws.Range("N70").Value = -5539.75
$ws.Range("H73").Value = 5399.8
$ws.Range("J73").Value = 4999.75
$ws.Range("L73").Value = 4999.75
$ws.Range("N73").Value = -6871.75
$ws.Range("H76").Value = 19665.666
$ws.Range("J76").Value = 19665.666
$ws.Range("L76").Value = 19665.666
$ws.Range("N76").Value = -20295.666
$ws.Range("H79").Value = 19665.666
$ws.Range("J79").Value = 19665.666
$ws.Range("L79").Value = 19665.666
$ws.Range("N79").Value = -21849.666
$ws.Range("H132").Value = 837994.4399999999
$ws.Range("I132").Value = 1166843.2
$ws.Range("J132").Value = 3224.2307
$ws.Range("K132").Value = 3500529.6
$ws.Range("L132").Value = 9672.6921
$ws.Range("M132").Value = -3497999.6
$ws.Range("N132").Value = -14732.6921
$ws.Range("H138").Value = 30000
$ws.Range("J138").Value = 30000
$ws.Range("L138").Value = 30000
$ws.Range("N138").Value = -40280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2079.457
$ws.Range("I132").Value = 826.4167
$ws.Range("K132").Value = 2479.2501
$ws.Range("M132").Value = 50.7498999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 42495.195
$ws.Range("I122").Value = 50386.617
$ws.Range("K122").Value = 151159.851
$ws.Range("M122").Value = -148709.851
